$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "IA1" batch (row 26) to "IA" ---
$ws.Range("B26").Value = "IA"

# --- Add J26 (empty, same look as the rest of the row) ---
$ws.Range("I26:I26").Copy()
$ws.Range("J26").PasteSpecial(-4122)
$ws.Range("J26").ClearContents()

# --- New row 27: a fresh "EC" batch / course 1 line, same layout as row 26 ---
$ws.Range("A26:I26").Copy()
$ws.Range("A27:I27").PasteSpecial(-4122)
$ws.Range("H23").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("J26").Copy()
$ws.Range("J27").PasteSpecial(-4122)

$ws.Range("A27").Value = "SE"
$ws.Range("B27").Value = "EC"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = "MAE101"
$ws.Range("F27").Value = "CSI101"
$ws.Range("G27").Value = "CEA201"
$ws.Range("H27").Value = "PRF192"
$ws.Range("I27").Value = "PHY101"
$ws.Range("J27").ClearContents()

# --- Three trailing blank-but-formatted rows (28-30), columns E:J ---
$ws.Range("E26:J26").Copy()
$ws.Range("E28:J28").PasteSpecial(-4122)
$ws.Range("E28:J28").ClearContents()

$ws.Range("E26:J26").Copy()
$ws.Range("E29:J29").PasteSpecial(-4122)
$ws.Range("E29:J29").ClearContents()
$ws.Range("H23").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").ClearContents()

$ws.Range("E26:J26").Copy()
$ws.Range("E30:J30").PasteSpecial(-4122)
$ws.Range("E30:J30").ClearContents()

# --- View state: active cell / selection ---
$ws.Range("M21").Select()
